# "Lista vocabulario.xlsx" — Tema 7 "Clothes to describe" section touch-up.
#
# The author filled a few missing English/Spanish synonym cells in columns
# B/D, fixed the capitalisation of two entries ("Falda" -> "falda",
# "Shoes" -> "shoes"), and extended the CONCATENATE helper formula in
# column H from H151 down to H176 (it previously stopped at H151).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Fill in the synonym cells that were previously left blank ----------
$ws.Range("B155").Value = "pants"       # trousers / pantalones synonym
$ws.Range("B162").Value = "waistcoat"   # vest synonym
$ws.Range("D162").Value = "chaleco"     # cuello pico synonym (Spanish)
$ws.Range("B164").Value = "jumper"      # sweater synonym
$ws.Range("B165").Value = "sweater"     # jumper synonym
$ws.Range("C169").Value = "falda"       # was "Falda" -> corrected casing
$ws.Range("A170").Value = "shoes"       # was "Shoes" -> corrected casing

# --- Extend the shared "Tema - Topic" helper formula down to row 176 ----
# Originally H131:H151 carried CONCATENATE(E,"-",F); the fill/formula
# was never continued when rows 152-176 were added, so H152:H176 were
# empty. Re-create the same formula pattern for those rows.
$ws.Range("H152:H176").Formula = '=CONCATENATE(E152," ","-"," ",F152)'

# --- Restore the cursor position the author ended up on -----------------
$ws.Range("I172").Select()
